$d = $word.ActiveDocument

# Get the last paragraph in the document body (before sectPr)
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)

# Insert a new paragraph after the last paragraph
$newRange = $lastPara.Range.InsertParagraphAfter()

# Move to the newly inserted paragraph and set its style + text
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = "Heading1"
$newPara.Range.Text = "CLASS ITEM COIN"
